$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to be treated as text so numeric-looking values
# (e.g. "592.73") are not auto-converted into floating point numbers,
# matching the original inline-string ("t=inlineStr") cell contents.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.758.89'
$ws.Range("E2").Value = '  +6.38%  '
$ws.Range("D3").Value = '3.579.23'
$ws.Range("E3").Value = '  +5.38%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '592.73'
$ws.Range("E5").Value = '  +5.85%  '
$ws.Range("D6").Value = '192.29'
$ws.Range("E6").Value = '  +9.17%  '
$ws.Range("D7").Value = '0.643'
$ws.Range("E7").Value = '  +1.70%  '
$ws.Range("D8").Value = '3.574.78'
$ws.Range("E8").Value = '  +5.52%  '
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").Value = '0.185'
$ws.Range("E10").Value = '  +6.08%  '
$ws.Range("D11").Value = '0.662'
$ws.Range("E11").Value = '  +4.22%  '
$ws.Range("D12").Value = '58.26'
$ws.Range("E12").Value = '  +8.75%  '
$ws.Range("D13").Value = '0.0000294'
$ws.Range("E13").Value = '  +5.79%  '
$ws.Range("D14").Value = '9.74'
$ws.Range("E14").Value = '  +5.80%  '
$ws.Range("D15").Value = '4.139.41'
$ws.Range("E15").Value = '  +5.12%  '
$ws.Range("D16").Value = '19.33'
$ws.Range("E16").Value = '  +5.76%  '
$ws.Range("D17").Value = '3.565.17'
$ws.Range("E17").Value = '  +5.45%  '
$ws.Range("D18").Value = '69.573.66'
$ws.Range("E18").Value = '  +6.06%  '
$ws.Range("D19").Value = '12.48'
$ws.Range("E19").Value = '  +5.37%  '
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("E21").Value = '  +4.95%  '
$ws.Range("D22").Value = '500.40'
$ws.Range("E22").Value = '  +4.48%  '
$ws.Range("D23").Value = '5.49'
$ws.Range("E23").Value = '  +11.28%  '
$ws.Range("D24").Value = '17.19'
$ws.Range("E24").Value = '  +20.06%  '
$ws.Range("D25").Value = '4.45'
$ws.Range("E25").Value = '  +8.54%  '
$ws.Range("D26").Value = '91.23'
$ws.Range("E26").Value = '  +1.71%  '
$ws.Range("D27").Value = '3.08'
$ws.Range("E27").Value = '  +5.35%  '
$ws.Range("D28").Value = '11.20'
$ws.Range("E28").Value = '  +5.09%  '
$ws.Range("D29").Value = '9.36'
$ws.Range("E29").Value = '  +7.23%  '
$ws.Range("D30").Value = '32.20'
$ws.Range("E30").Value = '  +2.93%  '
$ws.Range("D31").Value = '7.52'
$ws.Range("E31").Value = '  +14.62%  '
$ws.Range("D32").Value = '12.15'
$ws.Range("E32").Value = '  +5.73%  '
$ws.Range("D33").Value = '612.37'
$ws.Range("E33").Value = '  +6.58%  '
$ws.Range("D34").Value = '65.42'
$ws.Range("E34").Value = '  +3.79%  '
$ws.Range("D35").Value = '0.116'
$ws.Range("E35").Value = '  +7.33%  '
$ws.Range("D36").Value = '0.0₃0838'
$ws.Range("E36").Value = '  +13.10%  '
$ws.Range("E37").Value = '  +4.60%  '
$ws.Range("D38").Value = '38.15'
$ws.Range("E38").Value = '  +6.48%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = '0.399'
$ws.Range("E40").Value = '  +6.52%  '
$ws.Range("D41").Value = '3.64'
$ws.Range("E41").Value = '  -1.15%  '
$ws.Range("D42").Value = '3.325.50'
$ws.Range("E42").Value = '  +7.67%  '
$ws.Range("D43").Value = '3.14'
$ws.Range("E43").Value = '  +11.90%  '
$ws.Range("D44").Value = '2.72'
$ws.Range("E44").Value = '  +11.78%  '
$ws.Range("D45").Value = '0.0443'
$ws.Range("E45").Value = '  +6.34%  '
$ws.Range("D46").Value = '2.92'
$ws.Range("E46").Value = '  +19.52%  '
$ws.Range("D47").Value = '3.34'
$ws.Range("E47").Value = '  +5.09%  '
$ws.Range("E48").Value = '  +2.38%  '
$ws.Range("D49").Value = '9.14'
$ws.Range("E49").Value = '  +8.43%  '
$ws.Range("E50").Value = '  +4.49%  '
$ws.Range("D51").Value = '0.997'
$ws.Range("E51").Value = '  -0.23%  '

# Restore the default cell style so no residual number-format styling
# is left behind on the price column.
$ws.Range("D2:D51").Style = "Normal"

